# Regenerate merged AHB files
# 1) Rename header cells in row 1: "_old" -> "_FV2310", "_new" -> "_FV2404"
# 2) Freeze the header row (row 1)
# 3) Convert the data range A1:U62 into an Excel Table ("Table1")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldHeaders = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

$newHeaders = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $oldHeaders[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $newHeaders[$i]
}

# Freeze the top row
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the used range into a table
$rng = $ws.Range("A1:U62")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes, $null)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$wb.Save()
